$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '27.508.52', '  -5.13%  ')
    ,@(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.838.46', '  -4.52%  ')
    ,@(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.001', '  -0.33%  ')
    ,@(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '312.52', '  -3.86%  ')
    ,@(6, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.000', '  -0.28%  ')
    ,@(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.4238', '  -7.67%  ')
    ,@(8, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.3618', '  -5.28%  ')
    ,@(9, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '43.56', '  -4.50%  ')
    ,@(10, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.07177', '  -7.44%  ')
    ,@(11, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.8962', '  -8.45%  ')
    ,@(12, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '20.59', '  -8.73%  ')
    ,@(13, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.837.10', '  -4.28%  ')
    ,@(14, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '6.576', '  -5.66%  ')
    ,@(15, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.309', '  -7.02%  ')
    ,@(16, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.06796', '  -2.75%  ')
    ,@(17, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.001', '  -0.49%  ')
    ,@(18, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '77.17', '  -9.01%  ')
    ,@(19, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000008913', '  -6.06%  ')
    ,@(20, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.000', '  -0.27%  ')
    ,@(21, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '15.29', '  -8.53%  ')
    ,@(22, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '27.495.78', '  -5.25%  ')
    ,@(23, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '4.918', '  -8.02%  ')
    ,@(24, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '10.71', '  -3.51%  ')
    ,@(25, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.044.93', '  -5.24%  ')
    ,@(26, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '2.040', '  -0.78%  ')
    ,@(27, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '151.25', '  -4.38%  ')
    ,@(28, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '18.15', '  -4.65%  ')
    ,@(29, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '5.307', '  -5.69%  ')
    ,@(30, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '110.93', '  -5.86%  ')
    ,@(31, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '1.721', '  -6.56%  ')
    ,@(32, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.08866', '  -4.85%  ')
    ,@(33, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.7726', '  -10.68%  ')
    ,@(34, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '4.464', '  -12.60%  ')
    ,@(35, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.852', '  -5.47%  ')
    ,@(36, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.076', '  -13.76%  ')
    ,@(37, 'Frax', 'https://coinranking.com/coin/KfWtaeV1W+frax-frax', '1.000', '  -0.28%  ')
    ,@(38, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.05387', '  -5.49%  ')
    ,@(39, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '1.097', '  -4.91%  ')
    ,@(40, 'MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.940', '  -4.94%  ')
    ,@(41, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01913', '  -6.86%  ')
    ,@(42, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.5043', '  -8.51%  ')
    ,@(43, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '6.781', '  -9.20%  ')
    ,@(44, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.1631', '  -7.22%  ')
    ,@(45, 'Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.06611', '  -4.75%  ')
    ,@(46, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '8.169', '  -12.66%  ')
    ,@(47, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '106.01', '  -4.52%  ')
    ,@(48, 'Decentraland', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana', '0.4691', '  -9.35%  ')
    ,@(49, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '10.17', '  -9.24%  ')
    ,@(50, 'PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '0.9998', '  -0.35%  ')
    ,@(51, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '1.640', '  -7.22%  ')
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").NumberFormat = "@"
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").NumberFormat = "@"
    $ws.Range("E$r").Value = $row[4]
}
